$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing row (row 8) that is no longer part of the table
$ws.Rows.Item(8).Delete()

# Header row: updated date label (columns B/C headers "Ballgorithm"/"ESPN" stay the same)
$ws.Range("A1").Value = "NBA, Friday 9th Feb 2024"
$ws.Range("B1").Value = "Ballgorithm"
$ws.Range("C1").Value = "ESPN"

# Row 2
$ws.Range("A2").Value = "Atlanta Hawks (22-29) vs Philadelphia 76ers (30-20)"
$ws.Range("B2").Value = "Philadelphia 76ers (65.38%)"
$ws.Range("C2").Value = "Philadelphia 76ers (62.4%)"

# Row 3
$ws.Range("A3").Value = "Washington Wizards (9-41) vs Boston Celtics (39-12)"
$ws.Range("B3").Value = "Boston Celtics (88.89%) "
$ws.Range("C3").Value = "Boston Celtics (94.6%) "

# Row 4
$ws.Range("A4").Value = "Houston Rockets (23-27) vs Toronto Raptors (18-33)"
$ws.Range("B4").Value = "Houston Rockets (66.67%)"
$ws.Range("C4").Value = "Houston Rockets (60.7%)"

# Row 5
$ws.Range("A5").Value = "Charlotte Hornets (10-40) vs Milwaukee Bucks (33-19)"
$ws.Range("B5").Value = "Milwaukee Bucks (77.78%)"
$ws.Range("C5").Value = "Milwaukee Bucks (88.8%)"

# Row 6
$ws.Range("A6").Value = "Denver Nuggets (36-16) vs Sacramento Kings (29-21)"
$ws.Range("B6").Value = "Denver Nuggets (84.00%)"
$ws.Range("C6").Value = "Sacramento Kings (52.6%)"

# Row 7
$ws.Range("A7").Value = "New Orleans Pelicans (30-21) vs Los Angeles Lakers (28-26)"
$ws.Range("B7").Value = "Los Angeles Lakers (66.67%)"
$ws.Range("C7").Value = "New Orleans Pelicans (64.7%)"

# Update selection to match the saved workbook state
$ws.Range("A2").Select()
